$d = $word.ActiveDocument

$replacements = @(
    @{old="926×3="; new="534×3="},
    @{old="175×5="; new="234×2="},
    @{old="541×4="; new="648×9="},
    @{old="766×8="; new="169×7="},
    @{old="948×8="; new="604×7="},
    @{old="554×4="; new="889×9="},
    @{old="897×7="; new="696×7="},
    @{old="812×4="; new="193×4="},
    @{old="118×6="; new="687×7="},
    @{old="767×3="; new="800×9="},
    @{old="852×5="; new="974×6="},
    @{old="904×8="; new="508×2="},
    @{old="247×7="; new="965×4="},
    @{old="203×8="; new="376×9="},
    @{old="742×4="; new="705×8="},
    @{old="995×9="; new="105×9="},
    @{old="707×3="; new="195×8="},
    @{old="405×5="; new="901×2="},
    @{old="146×6="; new="289×6="},
    @{old="417×5="; new="966×6="},
    @{old="760×7="; new="780×8="},
    @{old="650×4="; new="399×9="},
    @{old="864×3="; new="227×2="},
    @{old="994×4="; new="284×4="},
    @{old="693×7="; new="549×3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
